# Auto-generated Excel COM-interop edit script
# Applies scheduled-runner price/profit updates to the Leve tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 306.18182
$ws.Range("I2").Value = 145.5
$ws.Range("J2").Value = 499
$ws.Range("K2").Value = 145.5
$ws.Range("L2").Value = 499
$ws.Range("M2").Value = -32.5
$ws.Range("N2").Value = -725

$ws.Range("H15").Value = 1331.8036
$ws.Range("I15").Value = 1331.8036
$ws.Range("K15").Value = 3995.4108
$ws.Range("M15").Value = -3826.4108

$ws.Range("H40").Value = 7559.222
$ws.Range("J40").Value = 10339.333
$ws.Range("L40").Value = 10339.333
$ws.Range("N40").Value = -10689.333

$ws.Range("H101").Value = 2921.8572
$ws.Range("I101").Value = 1489.3334
$ws.Range("J101").Value = 3996.25
$ws.Range("K101").Value = 4468.0002
$ws.Range("L101").Value = 11988.75
$ws.Range("M101").Value = -2846.0002
$ws.Range("N101").Value = -15232.75

$ws.Range("H121").Value = 849.3333
$ws.Range("J121").Value = 849.3333
$ws.Range("L121").Value = 2547.9999
$ws.Range("N121").Value = -6041.9999

$ws.Range("H135").Value = 2679.5789
$ws.Range("I135").Value = 2778.4443
$ws.Range("K135").Value = 25005.9987
$ws.Range("M135").Value = -22470.9987

$ws.Range("H138").Value = 3451.5972
$ws.Range("J138").Value = 3569.3281
$ws.Range("L138").Value = 10707.9843
$ws.Range("N138").Value = -20987.9843

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5786.4385
$ws.Range("I32").Value = 5172.3613
$ws.Range("K32").Value = 5172.3613
$ws.Range("M32").Value = -4885.3613

$ws.Range("H74").Value = 3029.923
$ws.Range("I74").Value = 2746.0908
$ws.Range("J74").Value = 4591
$ws.Range("K74").Value = 2746.0908
$ws.Range("L74").Value = 4591
$ws.Range("M74").Value = -1872.0908
$ws.Range("N74").Value = -6339

$ws.Range("H77").Value = 3029.923
$ws.Range("I77").Value = 2746.0908
$ws.Range("J77").Value = 4591
$ws.Range("K77").Value = 13730.454
$ws.Range("L77").Value = 22955
$ws.Range("M77").Value = -9362.454
$ws.Range("N77").Value = -31691

$ws.Range("H108").Value = 87385.8
$ws.Range("J108").Value = 87385.8
$ws.Range("L108").Value = 87385.8
$ws.Range("N108").Value = -95065.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2010
$ws.Range("I99").Value = 1935.6666
$ws.Range("J99").Value = 2307.3333
$ws.Range("K99").Value = 1935.6666
$ws.Range("L99").Value = 2307.3333
$ws.Range("M99").Value = -437.6666
$ws.Range("N99").Value = -5303.3333

$ws.Range("H105").Value = 11869.972
$ws.Range("I105").Value = 10277.4
$ws.Range("J105").Value = 15851.4
$ws.Range("K105").Value = 10277.4
$ws.Range("L105").Value = 15851.4
$ws.Range("M105").Value = -8530.4
$ws.Range("N105").Value = -19345.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 73390.266
$ws.Range("I31").Value = 3233
$ws.Range("K31").Value = 3233
$ws.Range("M31").Value = -2938

$ws.Range("H34").Value = 73390.266
$ws.Range("I34").Value = 3233
$ws.Range("K34").Value = 3233
$ws.Range("M34").Value = -3031

$ws.Range("H37").Value = 42331.332
$ws.Range("J37").Value = 47497
$ws.Range("L37").Value = 47497
$ws.Range("N37").Value = -47711

$ws.Range("H58").Value = 3728.9167
$ws.Range("I58").Value = 1701.5
$ws.Range("J58").Value = 7783.75
$ws.Range("K58").Value = 1701.5
$ws.Range("L58").Value = 7783.75
$ws.Range("M58").Value = -1498.5
$ws.Range("N58").Value = -8189.75

$ws.Range("H122").Value = 4193.2915
$ws.Range("I122").Value = 1653.3334
$ws.Range("J122").Value = 6733.25
$ws.Range("K122").Value = 4960.0002
$ws.Range("L122").Value = 20199.75
$ws.Range("M122").Value = -2510.0002
$ws.Range("N122").Value = -25099.75

$ws.Range("H136").Value = 3728.9167
$ws.Range("I136").Value = 1701.5
$ws.Range("J136").Value = 7783.75
$ws.Range("K136").Value = 5104.5
$ws.Range("L136").Value = 23351.25
$ws.Range("M136").Value = -2554.5
$ws.Range("N136").Value = -28451.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 5001.5
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5001.5
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 15004.5
$ws.Range("N31").Value = -15580.5
$ws.Range("M31").ClearContents()

$ws.Range("H44").Value = 103
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H70").Value = 14835.667
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 30000
$ws.Range("M70").Value = -29685

$ws.Range("H73").Value = 14835.667
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 30000
$ws.Range("M73").Value = -28908

$ws.Range("H76").Value = 12507.5
$ws.Range("I76").Value = 6000
$ws.Range("K76").Value = 18000
$ws.Range("M76").Value = -17617

$ws.Range("H79").Value = 12507.5
$ws.Range("I79").Value = 6000
$ws.Range("K79").Value = 18000
$ws.Range("M79").Value = -16674

$ws.Range("H113").Value = 1862.25
$ws.Range("J113").Value = 2379.6
$ws.Range("L113").Value = 7138.799999999999
$ws.Range("N113").Value = -11478.8

$ws.Range("H121").Value = 905.2857
$ws.Range("I121").Value = 898
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 2694
$ws.Range("L121").Value = 3000
$ws.Range("M121").Value = -1384
$ws.Range("N121").Value = -5620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 425.125
$ws.Range("I2").Value = 44.6
$ws.Range("K2").Value = 44.6
$ws.Range("M2").Value = 68.40000000000001

$ws.Range("H80").Value = 6115.091
$ws.Range("I80").Value = 2749.5
$ws.Range("J80").Value = 6863
$ws.Range("K80").Value = 2749.5
$ws.Range("L80").Value = 6863
$ws.Range("M80").Value = -1751.5
$ws.Range("N80").Value = -8859

$ws.Range("H83").Value = 6115.091
$ws.Range("I83").Value = 2749.5
$ws.Range("J83").Value = 6863
$ws.Range("K83").Value = 13747.5
$ws.Range("L83").Value = 34315
$ws.Range("M83").Value = -8755.5
$ws.Range("N83").Value = -44299

$ws.Range("H93").Value = 35275.7
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 35275.7
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 35275.7
$ws.Range("N93").Value = -39019.7
$ws.Range("M93").ClearContents()

$ws.Range("H109").Value = 64284
$ws.Range("J109").Value = 64284
$ws.Range("L109").Value = 64284
$ws.Range("N109").Value = -66364

$ws.Range("H138").Value = 78624.75
$ws.Range("J138").Value = 78624.75
$ws.Range("L138").Value = 78624.75
$ws.Range("N138").Value = -88904.75

$ws.Range("H139").Value = 93442
$ws.Range("J139").Value = 93442
$ws.Range("L139").Value = 93442
$ws.Range("N139").Value = -103722

$ws.Range("H141").Value = 84597.39999999999
$ws.Range("J141").Value = 84597.39999999999
$ws.Range("L141").Value = 84597.39999999999
$ws.Range("N141").Value = -94957.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 415.75
$ws.Range("I32").Value = 415.75
$ws.Range("K32").Value = 415.75
$ws.Range("M32").Value = -98.75

$ws.Range("H93").Value = 3281.5625
$ws.Range("I93").Value = 3291.4
$ws.Range("K93").Value = 3291.4
$ws.Range("M93").Value = -2043.4

$ws.Range("H136").Value = 3115.25
$ws.Range("I136").Value = 2034.4865
$ws.Range("J136").Value = 8827.857
$ws.Range("K136").Value = 6103.4595
$ws.Range("L136").Value = 26483.571
$ws.Range("M136").Value = -3553.4595
$ws.Range("N136").Value = -31583.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H12").Value = 9000
$ws.Range("J12").Value = 9000
$ws.Range("L12").Value = 9000
$ws.Range("N12").Value = -9284

$ws.Range("H107").Value = 1105.2368
$ws.Range("I107").Value = 1039.7778
$ws.Range("J107").Value = 1265.909
$ws.Range("K107").Value = 3119.3334
$ws.Range("L107").Value = 3797.727
$ws.Range("M107").Value = -1199.3334
$ws.Range("N107").Value = -7637.727000000001

$ws.Range("H113").Value = 497.2
$ws.Range("I113").Value = 488.17648
$ws.Range("K113").Value = 1464.52944
$ws.Range("M113").Value = 705.47056

$ws.Range("H132").Value = 1952.079
$ws.Range("J132").Value = 7671
$ws.Range("L132").Value = 23013
$ws.Range("N132").Value = -28073

$ws.Range("H136").Value = 2117.4707
$ws.Range("I136").Value = 1134.5217
$ws.Range("K136").Value = 3403.5651
$ws.Range("M136").Value = -853.5650999999998

